# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds strikeouts (K). The values below were recomputed from the
# underlying box-score data (K) to replace the previous "Strike#" values
# and are written directly into the sheet.
$kVals = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    14 = 1
}

foreach ($row in $kVals.Keys) {
    $ws.Range("G$row").Value = $kVals[$row]
}
